$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list figures (prices / 1h volume %) per the latest scrape.
# A leading apostrophe forces Excel to store the value as text (matching the
# original inline-string cells) while keeping the cell's General number format.
$ws.Range("D2").Value = '''42.874.87'
$ws.Range("E2").Value = '''  +0.12%  '
$ws.Range("D3").Value = '''2.534.42'
$ws.Range("E3").Value = '''  -0.89%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '''  +0.01%  '
$ws.Range("D5").Value = '''312.34'
$ws.Range("D6").Value = '''100.86'
$ws.Range("E6").Value = '''  +2.36%  '
$ws.Range("E7").Value = '''  -0.85%  '
$ws.Range("E8").Value = '''  +0.11%  '
$ws.Range("E9").Value = '''  -1.40%  '
$ws.Range("D10").Value = '''35.84'
$ws.Range("E10").Value = '''  +0.57%  '
$ws.Range("E11").Value = '''  -0.16%  '
$ws.Range("D12").Value = '''7.34'
$ws.Range("E12").Value = '''  -1.33%  '
$ws.Range("E13").Value = '''  +1.44%  '
$ws.Range("D14").Value = '''2.923.83'
$ws.Range("E14").Value = '''  -0.93%  '
$ws.Range("D15").Value = '''15.46'
$ws.Range("E15").Value = '''  -2.87%  '
$ws.Range("D16").Value = '''2.538.18'
$ws.Range("E16").Value = '''  -3.11%  '
$ws.Range("D17").Value = '''0.818'
$ws.Range("E17").Value = '''  -2.57%  '
$ws.Range("D18").Value = '''42.855.01'
$ws.Range("E18").Value = '''  +0.07%  '
$ws.Range("D19").Value = '''6.69'
$ws.Range("E19").Value = '''  -1.06%  '
$ws.Range("D20").Value = '''12.37'
$ws.Range("E20").Value = '''  +0.22%  '
$ws.Range("D22").Value = '''69.89'
$ws.Range("E22").Value = '''  +0.65%  '
$ws.Range("D23").Value = '''244.38'
$ws.Range("E23").Value = '''  -1.46%  '
$ws.Range("E24").Value = '''  -1.08%  '
$ws.Range("E25").Value = '''  -0.42%  '
$ws.Range("E26").Value = '''  +0.05%  '
$ws.Range("E27").Value = '''  -4.80%  '
$ws.Range("E28").Value = '''  -1.65%  '
$ws.Range("D29").Value = '''10.24'
$ws.Range("E29").Value = '''  +0.83%  '
$ws.Range("D30").Value = '''38.95'
$ws.Range("E30").Value = '''  -2.36%  '
$ws.Range("D31").Value = '''161.23'
$ws.Range("E31").Value = '''  +1.89%  '
$ws.Range("D32").Value = '''5.85'
$ws.Range("E32").Value = '''  +1.28%  '
$ws.Range("E33").Value = '''  +7.97%  '
$ws.Range("B34").Value = '''Hedera'
$ws.Range("C34").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.0792'
$ws.Range("E34").Value = '''  -0.68%  '
$ws.Range("B35").Value = '''WEMIXToken'
$ws.Range("C35").Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''2.66'
$ws.Range("E35").Value = '''  -0.21%  '
$ws.Range("E36").Value = '''  -1.78%  '
$ws.Range("D37").Value = '''3.13'
$ws.Range("E37").Value = '''  -5.12%  '
$ws.Range("D38").Value = '''1.98'
$ws.Range("E38").Value = '''  -5.55%  '
$ws.Range("E39").Value = '''  +0.01%  '
$ws.Range("E40").Value = '''  +0.08%  '
$ws.Range("D41").Value = '''4.18'
$ws.Range("E41").Value = '''  +2.19%  '
$ws.Range("D42").Value = '''22.05'
$ws.Range("E42").Value = '''  -3.72%  '
$ws.Range("E43").Value = '''  +4.23%  '
$ws.Range("E44").Value = '''  +0.20%  '
$ws.Range("D46").Value = '''2.004.07'
$ws.Range("E46").Value = '''  +0.66%  '
$ws.Range("D47").Value = '''9.25'
$ws.Range("E47").Value = '''  +3.16%  '
$ws.Range("D48").Value = '''2.776.51'
$ws.Range("E48").Value = '''  -1.06%  '
$ws.Range("D49").Value = '''0.192'
$ws.Range("E49").Value = '''  -0.47%  '
$ws.Range("D50").Value = '''79.74'
$ws.Range("E50").Value = '''  -2.00%  '
$ws.Range("D51").Value = '''72.60'
$ws.Range("E51").Value = '''  -1.44%  '
